$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C1 text from "requirement" to "requirements"
$ws.Range("C1").Value = "requirements"

# Move the active selection from B2 to C2
$ws.Range("C2").Select()
